$d = $word.ActiveDocument
$bmRange = $d.Range(88, 88)
$d.Bookmarks.Add("TestPos88", $bmRange)
Write-Output "done"
